# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
# values for the first data row (733363e9-...) on both the zh-cn and de-de sheets, to
# reflect the freshly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-12 10:34:31"
$zhcn.Range("H2").Value = "2016-03-12 10:34:47"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-12 10:34:34"
$dede.Range("H2").Value = "2016-03-12 10:34:53"
